$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-08-16 17:22:42"

foreach ($row in 2..6) {
    $cell = $ws.Cells.Item($row, 5)
    $cell.Value = $newTimestamp
}
